$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.904.06"
$ws.Range("E2").Value = "  -2.96%  "
$ws.Range("D3").Value = "1.625.59"
$ws.Range("E3").Value = "  -2.90%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").Value = "'1.007"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").Value = "'307.60"
$ws.Range("E6").Value = "  -2.20%  "
$ws.Range("D7").Value = "'0.3904"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "'0.3824"
$ws.Range("E8").Value = "  -2.80%  "
$ws.Range("D9").Value = "'1.009"
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("D10").Value = "'49.78"
$ws.Range("E10").Value = "  -4.36%  "
$ws.Range("D11").Value = "'1.346"
$ws.Range("E11").Value = "  -3.01%  "
$ws.Range("D12").Value = "'0.08483"
$ws.Range("E12").Value = "  -1.74%  "
$ws.Range("D13").Value = "'23.57"
$ws.Range("E13").Value = "  -6.09%  "
$ws.Range("D14").Value = "'6.983"
$ws.Range("E14").Value = "  -4.34%  "
$ws.Range("D15").Value = "'0.00001271"
$ws.Range("E15").Value = "  -3.11%  "
$ws.Range("D16").Value = "'7.415"
$ws.Range("E16").Value = "  -4.11%  "
$ws.Range("D17").Value = "1.635.89"
$ws.Range("E17").Value = "  -2.24%  "
$ws.Range("D18").Value = "'92.90"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("D19").Value = "'0.06910"
$ws.Range("E19").Value = "  -2.05%  "
$ws.Range("D20").Value = "'19.96"
$ws.Range("E20").Value = "  -2.31%  "
$ws.Range("D21").Value = "'6.839"
$ws.Range("E21").Value = "  -3.04%  "
$ws.Range("D22").Value = "'1.006"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").Value = "'13.40"
$ws.Range("E23").Value = "  -3.85%  "
$ws.Range("D24").Value = "23.914.28"
$ws.Range("E24").Value = "  -2.96%  "
$ws.Range("D25").Value = "'2.378"
$ws.Range("E25").Value = "  +1.12%  "
$ws.Range("D26").Value = "'2.814"
$ws.Range("E26").Value = "  +3.05%  "
$ws.Range("D27").Value = "'22.03"
$ws.Range("E27").Value = "  -5.15%  "
$ws.Range("D28").Value = "'157.33"
$ws.Range("E28").Value = "  -3.25%  "
$ws.Range("D29").Value = "'138.65"
$ws.Range("E29").Value = "  -5.44%  "
$ws.Range("D30").Value = "'5.236"
$ws.Range("E30").Value = "  -8.47%  "
$ws.Range("D31").Value = "'7.780"
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("D32").Value = "'2.438"
$ws.Range("E32").Value = "  -0.89%  "
$ws.Range("D33").Value = "1.815.81"
$ws.Range("E33").Value = "  -2.35%  "
$ws.Range("D34").Value = "'0.07932"
$ws.Range("E34").Value = "  -5.27%  "
$ws.Range("D35").Value = "'6.636"
$ws.Range("E35").Value = "  -3.98%  "
$ws.Range("D36").Value = "'0.02875"
$ws.Range("E36").Value = "  -5.17%  "
$ws.Range("D37").Value = "'0.9436"
$ws.Range("E37").Value = "  -4.23%  "
$ws.Range("D38").Value = "'0.2649"
$ws.Range("E38").Value = "  -5.15%  "
$ws.Range("D39").Value = "'0.09160"
$ws.Range("E39").Value = "  -3.03%  "
$ws.Range("D40").Value = "'10.21"
$ws.Range("E40").Value = "  -3.07%  "
$ws.Range("D41").Value = "'1.407"
$ws.Range("E41").Value = "  -9.12%  "
$ws.Range("D42").Value = "'0.7419"
$ws.Range("E42").Value = "  -6.13%  "
$ws.Range("D43").Value = "'13.00"
$ws.Range("E43").Value = "  -3.49%  "
$ws.Range("D44").Value = "'15.89"
$ws.Range("E44").Value = "  -3.42%  "
$ws.Range("D45").Value = "'0.6808"
$ws.Range("E45").Value = "  -4.22%  "
$ws.Range("D46").Value = "'2.428"
$ws.Range("E46").Value = "  -5.16%  "
$ws.Range("D47").Value = "'4.073"
$ws.Range("E47").Value = "  -2.92%  "
$ws.Range("D48").Value = "'1.007"
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("D49").Value = "'0.08247"
$ws.Range("D50").Value = "'132.31"
$ws.Range("E50").Value = "  -3.61%  "
$ws.Range("D51").Value = "'1.243"
$ws.Range("E51").Value = "  -6.02%  "
